# STAI updated Sep 10
# The "State-Trait Anxiety Inventory" row's abbreviation/block-name become
# the "-Y2" variant, and the "State-Trait Anxiety Inventory - Y - Present"
# row's abbreviation/block-name become the "-Y1" variant.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: State-Trait Anxiety Inventory – Y – Present
$ws.Range("B27").Value = "STAI-Y1"

# Row 26: State-Trait Anxiety Inventory
$ws.Range("B26").Value = "STAI-Y2"
$ws.Range("C26").Value = "STAIY2_likert_block"

$ws.Range("C27").Value = "STAIY1_likert_block"

# Move selection to D29 (matches the saved view state in the target file)
$ws.Range("D29").Select()
